$wb = $excel.ActiveWorkbook

# ALC row 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 861.6
$ws.Range("I4").Value = 379.14285
$ws.Range("K4").Value = 379.14285
$ws.Range("M4").Value = -265.14285

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 10885
$ws.Range("I19").Value = 1998.1666
$ws.Range("J19").Value = 14217.5625
$ws.Range("K19").Value = 1998.1666
$ws.Range("L19").Value = 14217.5625
$ws.Range("M19").Value = -1823.1666
$ws.Range("N19").Value = -14567.5625

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 944.41174
$ws.Range("I88").Value = 942.125
$ws.Range("J88").Value = 946.44446
$ws.Range("K88").Value = 942.125
$ws.Range("L88").Value = 946.44446
$ws.Range("M88").Value = -536.125
$ws.Range("N88").Value = -1758.44446

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 944.41174
$ws.Range("I91").Value = 942.125
$ws.Range("J91").Value = 946.44446
$ws.Range("K91").Value = 942.125
$ws.Range("L91").Value = 946.44446
$ws.Range("M91").Value = 461.875
$ws.Range("N91").Value = -3754.44446

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 9420.714
$ws.Range("I100").Value = 10157.5
$ws.Range("K100").Value = 10157.5
$ws.Range("M100").Value = -9616.5

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3211
$ws.Range("I132").Value = 3318.182
$ws.Range("J132").Value = 2503.6
$ws.Range("K132").Value = 9954.545999999998
$ws.Range("L132").Value = 7510.799999999999
$ws.Range("M132").Value = -7424.545999999998
$ws.Range("N132").Value = -12570.8

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 90910650
$ws.Range("I137").Value = 100001210
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 300003630
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -300001080
$ws.Range("N137").Value = -20100

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1525.228
$ws.Range("I32").Value = 1593.963
$ws.Range("J32").Value = 288
$ws.Range("K32").Value = 1593.963
$ws.Range("L32").Value = 288
$ws.Range("M32").Value = -1306.963
$ws.Range("N32").Value = -862

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1865.4546
$ws.Range("I74").Value = 1112.6842
$ws.Range("K74").Value = 1112.6842
$ws.Range("M74").Value = -238.6841999999999

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1865.4546
$ws.Range("I77").Value = 1112.6842
$ws.Range("K77").Value = 5563.420999999999
$ws.Range("M77").Value = -1195.420999999999

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 873.4
$ws.Range("I97").Value = 737.0833
$ws.Range("K97").Value = 737.0833
$ws.Range("M97").Value = -241.0833

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 6114513.5
$ws.Range("I102").Value = 7576474.5
$ws.Range("J102").Value = 266670
$ws.Range("K102").Value = 7576474.5
$ws.Range("L102").Value = 266670
$ws.Range("M102").Value = -7574852.5
$ws.Range("N102").Value = -269914

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 38463644
$ws.Range("I132").Value = 41668616
$ws.Range("K132").Value = 125005848
$ws.Range("M132").Value = -125003318

# BSM row 70
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 300000
$ws.Range("J70").Value = 300000
$ws.Range("L70").Value = 300000
$ws.Range("N70").Value = -300586

# BSM row 73
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H73").Value = 300000
$ws.Range("J73").Value = 300000
$ws.Range("L73").Value = 300000
$ws.Range("N73").Value = -302028

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17859312
$ws.Range("I86").Value = 29414348
$ws.Range("J86").Value = 1529.7273
$ws.Range("K86").Value = 29414348
$ws.Range("L86").Value = 1529.7273
$ws.Range("M86").Value = -29413225
$ws.Range("N86").Value = -3775.7273

# BSM row 88
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 45000
$ws.Range("J88").Value = 45000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45812

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 17859312
$ws.Range("I89").Value = 29414348
$ws.Range("J89").Value = 1529.7273
$ws.Range("K89").Value = 147071740
$ws.Range("L89").Value = 7648.636500000001
$ws.Range("M89").Value = -147066124
$ws.Range("N89").Value = -18880.6365

# BSM row 91
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 45000
$ws.Range("J91").Value = 45000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47808

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3571.5
$ws.Range("I94").Value = 3661.9412
$ws.Range("J94").Value = 3264
$ws.Range("K94").Value = 3661.9412
$ws.Range("L94").Value = 3264
$ws.Range("M94").Value = -3210.9412
$ws.Range("N94").Value = -4166

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1087.1177
$ws.Range("I99").Value = 1055.6
$ws.Range("J99").Value = 1132.1428
$ws.Range("K99").Value = 1055.6
$ws.Range("L99").Value = 1132.1428
$ws.Range("M99").Value = 442.4000000000001
$ws.Range("N99").Value = -4128.1428

# BSM row 102
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 20278
$ws.Range("I102").Value = 20278
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 20278
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -17033

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1916.3334
$ws.Range("I105").Value = 1905.3077
$ws.Range("J105").Value = 1934.25
$ws.Range("K105").Value = 1905.3077
$ws.Range("L105").Value = 1934.25
$ws.Range("M105").Value = -158.3077000000001
$ws.Range("N105").Value = -5428.25

# BSM row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 59999
$ws.Range("J130").Value = 59999
$ws.Range("L130").Value = 59999
$ws.Range("N130").Value = -70039

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2720.4666
$ws.Range("I99").Value = 2280.7
$ws.Range("J99").Value = 3600
$ws.Range("K99").Value = 2280.7
$ws.Range("L99").Value = 3600
$ws.Range("M99").Value = -782.6999999999998
$ws.Range("N99").Value = -6596

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2510.2942
$ws.Range("I122").Value = 2131.9
$ws.Range("J122").Value = 3050.8572
$ws.Range("K122").Value = 6395.700000000001
$ws.Range("L122").Value = 9152.571599999999
$ws.Range("M122").Value = -3945.700000000001
$ws.Range("N122").Value = -14052.5716

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2720.4666
$ws.Range("I126").Value = 2280.7
$ws.Range("J126").Value = 3600
$ws.Range("K126").Value = 6842.099999999999
$ws.Range("L126").Value = 10800
$ws.Range("M126").Value = -4372.099999999999
$ws.Range("N126").Value = -15740

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 248
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 113448.45
$ws.Range("I121").Value = 1012
$ws.Range("K121").Value = 3036
$ws.Range("M121").Value = -1726

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1643.0769
$ws.Range("I107").Value = 508.85715
$ws.Range("K107").Value = 508.85715
$ws.Range("M107").Value = 1411.14285

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 9742.9375
$ws.Range("I126").Value = 11989.417
$ws.Range("J126").Value = 3003.5
$ws.Range("K126").Value = 35968.251
$ws.Range("L126").Value = 9010.5
$ws.Range("M126").Value = -33498.251
$ws.Range("N126").Value = -13950.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2728.4048
$ws.Range("I132").Value = 1793.5161
$ws.Range("K132").Value = 5380.5483
$ws.Range("M132").Value = -2850.5483

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 150000
$ws.Range("J134").Value = 150000
$ws.Range("L134").Value = 450000
$ws.Range("N134").Value = -455070

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1662.1818
$ws.Range("I7").Value = 1584.875
$ws.Range("K7").Value = 1584.875
$ws.Range("M7").Value = -1472.875

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1423
$ws.Range("I16").Value = 1341.0769
$ws.Range("K16").Value = 1341.0769
$ws.Range("M16").Value = -1171.0769

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2680.25
$ws.Range("I40").Value = 2682.6667
$ws.Range("K40").Value = 2682.6667
$ws.Range("M40").Value = -2546.6667

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2228.4666
$ws.Range("I46").Value = 989
$ws.Range("J46").Value = 2538.3333
$ws.Range("K46").Value = 989
$ws.Range("L46").Value = 2538.3333
$ws.Range("M46").Value = -801
$ws.Range("N46").Value = -2914.3333

# LTW row 57
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 40000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1265.4286
$ws.Range("I93").Value = 1262.4615
$ws.Range("K93").Value = 1262.4615
$ws.Range("M93").Value = -14.46149999999989

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3166.6667
$ws.Range("I100").Value = 2750
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 2750
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2209
$ws.Range("N100").Value = -5082

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3140.3704
$ws.Range("I122").Value = 2642.8667
$ws.Range("K122").Value = 7928.6001
$ws.Range("M122").Value = -5478.6001

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1662.1818
$ws.Range("I126").Value = 1584.875
$ws.Range("K126").Value = 4754.625
$ws.Range("M126").Value = -2284.625

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1896.8
$ws.Range("I100").Value = 2143.5518
$ws.Range("K100").Value = 4287.1036
$ws.Range("M100").Value = -3746.1036

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 727.1429000000001
$ws.Range("I107").Value = 698.3333
$ws.Range("K107").Value = 2094.9999
$ws.Range("M107").Value = -174.9998999999998

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2312.2
$ws.Range("I122").Value = 2258
$ws.Range("K122").Value = 6774
$ws.Range("M122").Value = -4324

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1600
$ws.Range("I126").Value = 1576.4706
$ws.Range("K126").Value = 4729.4118
$ws.Range("M126").Value = -2259.4118
